$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records to append: date (Excel serial), nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44301, 0, 7, 178.0264496439471),
    @(44302, 0, 6, 152.5940996948118),
    @(44303, 0, 3, 76.2970498474059)
)

$startRow = 227
$lastRow = $startRow - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy formatting (style) of the date cell from the previous row before writing the new value
    $ws.Cells.Item($lastRow, 1).Copy($ws.Cells.Item($r, 1))

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $lastRow = $r
}
